# Connor_Readnour_QA_Resume.docx edit:
#  1) Split the "ARCHIBUS IWMS software" run so the new phrase
#     "Web Central and Smart Client " is inserted before "IWMS", as three
#     separate runs (matching how Word splits a run when new text with the
#     same formatting is typed into the middle of existing text).
#  2) The "_GoBack" bookmark (which Word stamps at the location of the most
#     recent edit) moves from its old spot (mid-sentence in the next bullet,
#     "Developed over 25 r|esponsive...") to the empty paragraph that
#     immediately follows the paragraph that was just edited.

$d = $word.ActiveDocument

$wNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Step 1: find + rewrite the edited bullet paragraph.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$editIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Created and modified hundreds of views*") {
        $editIdx = $i
        break
    }
}
if ($editIdx -eq -1) { throw "Could not find target paragraph (Created and modified hundreds of views...)" }

$editPara = $paras.Item($editIdx)
$editRange = $editPara.Range

# Shared run formatting used throughout this bullet (unchanged by the edit).
$rPr = '<w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'

$part1 = "• Created and modified hundreds of views, home pages and reports in the ARCHIBUS "
$part2 = "Web Central and Smart Client IWMS software "
$part3 = "based on customer specifications."

$newParaXml = '<w:p ' + $wNS + ' w:rsidR="005A5603" w:rsidRPr="005E64E9" w:rsidRDefault="005A5603" w:rsidP="005A5603">' +
    '<w:pPr>' + $rPr + '</w:pPr>' +
    '<w:r w:rsidRPr="005E64E9">' + $rPr + '<w:t xml:space="preserve">' + $part1 + '</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $part2 + '</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>' + $part3 + '</w:t></w:r>' +
    '</w:p>'

$editRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# Step 2: add the "_GoBack" bookmark to the (empty) paragraph right after
# the one we just edited.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$afterPara = $paras.Item($editIdx + 1)
$afterRange = $afterPara.Range

$afterParaXml = '<w:p ' + $wNS + ' w:rsidR="005A5603" w:rsidRPr="005E64E9" w:rsidRDefault="005A5603" w:rsidP="005A5603">' +
    '<w:pPr>' + $rPr + '</w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$afterRange.InsertXML($afterParaXml)

# ---------------------------------------------------------------------
# Step 3: remove the old "_GoBack" bookmark from the "Developed over 25
# r..." paragraph (the bookmark used to sit mid-run there).
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$oldBmIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Developed over 25 r*") {
        $oldBmIdx = $i
        break
    }
}
if ($oldBmIdx -eq -1) { throw "Could not find paragraph that used to hold the _GoBack bookmark" }

$oldBmPara = $paras.Item($oldBmIdx)
$oldBmRange = $oldBmPara.Range

$rPrBoots = '<w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'

$oldBmParaXml = '<w:p ' + $wNS + ' w:rsidR="005D3E75" w:rsidRPr="00330DDF" w:rsidRDefault="00042101" w:rsidP="005D3E75">' +
    '<w:pPr>' + $rPrBoots + '</w:pPr>' +
    '<w:r>' + $rPrBoots + '<w:t>' + "• Developed over 25 r" + '</w:t></w:r>' +
    '<w:r w:rsidR="005D3E75" w:rsidRPr="00330DDF">' + $rPrBoots + '<w:t>esponsive SpaceView Boots</w:t></w:r>' +
    '<w:r w:rsidR="005D3E75">' + $rPrBoots + '<w:t xml:space="preserve">trap webpages for RSC customers, </w:t></w:r>' +
    '<w:r w:rsidR="005D3E75" w:rsidRPr="00330DDF">' + $rPrBoots + '<w:t>us</w:t></w:r>' +
    '<w:r w:rsidR="005D3E75">' + $rPrBoots + '<w:t>ing HTML5, CSS3, and JavaScript, as well as unit testing hundreds of design and functionality updates.</w:t></w:r>' +
    '</w:p>'

$oldBmRange.InsertXML($oldBmParaXml)

Write-Output "done"
